$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wilaya label on row 17 (drop the "10 - " prefix)
$ws.Range("C17").Value = "Bouira"

# Append new rows 18 and 19
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 23
$ws.Range("C18").Value = "Constantine"
$ws.Range("D18").Value = "Constantine"
$ws.Range("E18").Value = "666317ac09c6d4281f17aa37"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 23
$ws.Range("C19").Value = "27 - Mostaganem"
$ws.Range("D19").Value = "Ain Tedeles"
$ws.Range("E19").Value = "666317f909c6d4281f17aa3a"
